$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as exact text (preserve trailing zeros / formatting)
# by pre-formatting the specific target cells as Text before assigning their string values.
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D20","D21","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates per the diff
$ws.Range("D2").Value = '25.767.31'
$ws.Range("E2").Value = '  -2.72%  '
$ws.Range("D3").Value = '1.744.15'
$ws.Range("E3").Value = '  -5.12%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '238.40'
$ws.Range("E5").Value = '  -8.94%  '
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.5033'
$ws.Range("E7").Value = '  -6.70%  '
$ws.Range("D8").Value = '41.96'
$ws.Range("E8").Value = '  -6.46%  '
$ws.Range("D9").Value = '0.2734'
$ws.Range("E9").Value = '  -9.27%  '
$ws.Range("D10").Value = '0.06148'
$ws.Range("E10").Value = '  -10.96%  '
$ws.Range("D11").Value = '1.746.29'
$ws.Range("E11").Value = '  -4.98%  '
$ws.Range("D12").Value = '0.06920'
$ws.Range("E12").Value = '  -3.10%  '
$ws.Range("D13").Value = '15.49'
$ws.Range("E13").Value = '  -12.33%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '0.6042'
$ws.Range("E14").Value = '  -18.11%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '4.528'
$ws.Range("E15").Value = '  -9.41%  '
$ws.Range("D16").Value = '77.16'
$ws.Range("E16").Value = '  -13.36%  '
$ws.Range("D17").Value = '0.9991'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '25.775.21'
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("D20").Value = '0.000006875'
$ws.Range("E20").Value = '  -12.97%  '
$ws.Range("D21").Value = '11.60'
$ws.Range("E21").Value = '  -16.30%  '
$ws.Range("D22").Value = '1.967.57'
$ws.Range("E22").Value = '  -5.46%  '
$ws.Range("D23").Value = '4.059'
$ws.Range("E23").Value = '  -11.76%  '
$ws.Range("D24").Value = '5.246'
$ws.Range("E24").Value = '  -12.33%  '
$ws.Range("D25").Value = '8.182'
$ws.Range("E25").Value = '  -11.13%  '
$ws.Range("D26").Value = '137.88'
$ws.Range("E26").Value = '  -3.44%  '
$ws.Range("D27").Value = '1.459'
$ws.Range("E27").Value = '  -15.13%  '
$ws.Range("D28").Value = '1.820'
$ws.Range("E29").Value = '  -11.76%  '
$ws.Range("D30").Value = '103.83'
$ws.Range("E30").Value = '  -6.62%  '
$ws.Range("D31").Value = '0.08133'
$ws.Range("E31").Value = '  -8.18%  '
$ws.Range("D32").Value = '3.709'
$ws.Range("E32").Value = '  -12.71%  '
$ws.Range("D33").Value = '3.488'
$ws.Range("E33").Value = '  -13.87%  '
$ws.Range("D34").Value = '0.04555'
$ws.Range("E34").Value = '  -5.95%  '
$ws.Range("D35").Value = '0.9990'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").Value = '2.614'
$ws.Range("E36").Value = '  -10.63%  '
$ws.Range("D37").Value = '0.9852'
$ws.Range("D38").Value = '0.6111'
$ws.Range("E38").Value = '  -16.40%  '
$ws.Range("E39").Value = '  -13.57%  '
$ws.Range("D40").Value = '0.01555'
$ws.Range("E40").Value = '  -9.47%  '
$ws.Range("D41").Value = '1.932'
$ws.Range("E41").Value = '  -14.58%  '
$ws.Range("D42").Value = '0.9993'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '101.80'
$ws.Range("E43").Value = '  -5.83%  '
$ws.Range("D44").Value = '0.3855'
$ws.Range("E44").Value = '  -18.26%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.7343'
$ws.Range("E45").Value = '  -18.68%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '5.024'
$ws.Range("E46").Value = '  -14.87%  '
$ws.Range("D47").Value = '0.05371'
$ws.Range("E47").Value = '  -6.84%  '
$ws.Range("D48").Value = '0.1112'
$ws.Range("E48").Value = '  -11.36%  '
$ws.Range("D49").Value = '5.983'
$ws.Range("E49").Value = '  -19.14%  '
$ws.Range("D50").Value = '30.19'
$ws.Range("E50").Value = '  -13.27%  '
$ws.Range("D51").Value = '52.56'
$ws.Range("E51").Value = '  -12.55%  '
